$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.100.72"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "3.470.76"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.57"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.09"
$ws.Range("E6").Value = "  -2.30%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.481"
$ws.Range("E8").Value = "  -1.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.63"
$ws.Range("E9").Value = "  +5.35%  "
$ws.Range("E10").Value = "  -1.47%  "
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").Value = "4.062.96"
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.119"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("E14").Value = "  -2.77%  "
$ws.Range("D15").Value = "3.473.21"
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("D16").Value = "64.045.54"
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.25"
$ws.Range("E17").Value = "  -6.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.96"
$ws.Range("E18").Value = "  +0.60%  "
$ws.Range("E19").Value = "  -1.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.38"
$ws.Range("E20").Value = "  -1.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "383.61"
$ws.Range("E21").Value = "  -2.45%  "
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("D23").Value = "3.610.73"
$ws.Range("E23").Value = "  -0.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.81"
$ws.Range("E24").Value = "  +0.68%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.60"
$ws.Range("E26").Value = "  +0.94%  "
$ws.Range("E27").Value = "  -3.25%  "
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("E30").Value = "  -4.79%  "
$ws.Range("E31").Value = "  -4.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.93"
$ws.Range("E32").Value = "  -4.22%  "
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.152"
$ws.Range("E33").Value = "  +1.25%  "
$ws.Range("B34").Value = "RenzoRestakedETH"
$ws.Range("C34").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D34").Value = "3.497.98"
$ws.Range("E34").Value = "  -0.65%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "22.88"
$ws.Range("E36").Value = "  -2.39%  "
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.74"
$ws.Range("E38").Value = "  -2.47%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.50"
$ws.Range("E39").Value = "  -3.95%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "161.89"
$ws.Range("E40").Value = "  -1.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0777"
$ws.Range("E41").Value = "  -0.89%  "
$ws.Range("E42").Value = "  -1.15%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.23"
$ws.Range("E44").Value = "  -1.57%  "
$ws.Range("E45").Value = "  -2.72%  "
$ws.Range("E46").Value = "  -2.43%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.47"
$ws.Range("E47").Value = "  -7.30%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.13"
$ws.Range("E48").Value = "  -4.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.70"
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.900"
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("D51").Value = "2.323.86"
$ws.Range("E51").Value = "  -5.65%  "
